# Weekly update: insert a new Cilantro price record ahead of the existing
# history for "Terminal La Palmera de La Serena", pushing the old rows
# 146-152 down to 147-153.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 146:152 down by one row.
$ws.Rows("146:146").Insert()

# Populate the newly inserted row 146 with the new weekly observation.
$ws.Range("A146").Value = 8
$ws.Range("B146").Value = "Terminal La Palmera de La Serena"
$ws.Range("C146").Value = "Coquimbo"
$ws.Range("D146").Value = 44747
$ws.Range("E146").Value = 4
$ws.Range("F146").Value = 100112040
$ws.Range("G146").Value = "Cilantro"
$ws.Range("H146").Value = "Sin especificar"
$ws.Range("I146").Value = "Primera"
$ws.Range("J146").Value = 2400
$ws.Range("K146").Value = 1500
$ws.Range("L146").Value = 2000
$ws.Range("M146").Value = 1750
$ws.Range("N146").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O146").Value = "Provincia del Elquí"
$ws.Range("P146").Value = 1167
$ws.Range("Q146").Value = 1.5
$ws.Range("R146").Value = "Hortaliza"
